$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 57.38695133333332
$ws.Range("H2").Value = 172.160854
$ws.Range("I2").Value = 0.6542464432660997
$ws.Range("J2").Value = 0.6542464432660998
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.4526816666666666
$ws.Range("N2").Value = 1.358045
$ws.Range("O2").Value = 0.0009417166111561725
$ws.Range("P2").Value = 0.0009417166111561725
$ws.Range("Q2").Value = 25.97802077449222
$ws.Range("R2").Value = 233.8021869704299
$ws.Range("S2").Value = 0.0006161147434135306
$ws.Range("T2").Value = 0.0006161147434135306

# Row 3
$ws.Range("G3").Value = 57.38695133333332
$ws.Range("H3").Value = 172.160854
$ws.Range("I3").Value = 0.6542464432660997
$ws.Range("J3").Value = 0.6542464432660998
$ws.Range("M3").Value = 0.35382
$ws.Range("N3").Value = 1.06146
$ws.Range("O3").Value = 0.0007360540439218367
$ws.Range("P3").Value = 0.0007360540439218367
$ws.Range("Q3").Value = 20.30465112076
$ws.Range("R3").Value = 182.7418600868399
$ws.Range("S3").Value = 0.0004815607402874912
$ws.Range("T3").Value = 0.0004815607402874913

# Row 4
$ws.Range("G4").Value = 57.38695133333332
$ws.Range("H4").Value = 172.160854
$ws.Range("I4").Value = 0.6542464432660997
$ws.Range("J4").Value = 0.6542464432660998
$ws.Range("M4").Value = 479.8918966666667
$ws.Range("N4").Value = 1439.67569
$ws.Range("O4").Value = 0.998322229344922
$ws.Range("P4").Value = 0.998322229344922
$ws.Range("Q4").Value = 27539.53291927102
$ws.Range("R4").Value = 247855.7962734392
$ws.Range("S4").Value = 0.6531487677823987
$ws.Range("T4").Value = 0.6531487677823988

# Row 5
$ws.Range("I5").Value = 0.16357689713892
$ws.Range("J5").Value = 0.16357689713892
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.4526816666666666
$ws.Range("N5").Value = 1.358045
$ws.Range("O5").Value = 0.0009417166111561725
$ws.Range("P5").Value = 0.0009417166111561725
$ws.Range("Q5").Value = 6.495112164291111
$ws.Range("R5").Value = 58.45600947862
$ws.Range("S5").Value = 0.0001540430812371055
$ws.Range("T5").Value = 0.0001540430812371055

# Row 6
$ws.Range("I6").Value = 0.16357689713892
$ws.Range("J6").Value = 0.16357689713892
$ws.Range("M6").Value = 0.35382
$ws.Range("N6").Value = 1.06146
$ws.Range("O6").Value = 0.0007360540439218367
$ws.Range("P6").Value = 0.0007360540439218367
$ws.Range("Q6").Value = 5.07663719384
$ws.Range("R6").Value = 45.68973474456
$ws.Range("S6").Value = 0.0001204014366312884
$ws.Range("T6").Value = 0.0001204014366312884

# Row 7
$ws.Range("I7").Value = 0.16357689713892
$ws.Range("J7").Value = 0.16357689713892
$ws.Range("M7").Value = 479.8918966666667
$ws.Range("N7").Value = 1439.67569
$ws.Range("O7").Value = 0.998322229344922
$ws.Range("P7").Value = 0.998322229344922
$ws.Range("Q7").Value = 6885.526684869205
$ws.Range("R7").Value = 61969.74016382285
$ws.Range("S7").Value = 0.1633024526210516
$ws.Range("T7").Value = 0.1633024526210516

# Row 8
$ws.Range("G8").Value = 14.516389
$ws.Range("H8").Value = 43.549167
$ws.Range("I8").Value = 0.1654957381714162
$ws.Range("J8").Value = 0.1654957381714162
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.4526816666666666
$ws.Range("N8").Value = 1.358045
$ws.Range("O8").Value = 0.0009417166111561725
$ws.Range("P8").Value = 0.0009417166111561725
$ws.Range("Q8").Value = 6.571303166501666
$ws.Range("R8").Value = 59.14172849851499
$ws.Range("S8").Value = 0.0001558500857115753
$ws.Range("T8").Value = 0.0001558500857115753

# Row 9
$ws.Range("G9").Value = 14.516389
$ws.Range("H9").Value = 43.549167
$ws.Range("I9").Value = 0.1654957381714162
$ws.Range("J9").Value = 0.1654957381714162
$ws.Range("M9").Value = 0.35382
$ws.Range("N9").Value = 1.06146
$ws.Range("O9").Value = 0.0007360540439218367
$ws.Range("P9").Value = 0.0007360540439218367
$ws.Range("Q9").Value = 5.136188755979999
$ws.Range("R9").Value = 46.22569880381999
$ws.Range("S9").Value = 0.0001218138073329003
$ws.Range("T9").Value = 0.0001218138073329003

# Row 10
$ws.Range("G10").Value = 14.516389
$ws.Range("H10").Value = 43.549167
$ws.Range("I10").Value = 0.1654957381714162
$ws.Range("J10").Value = 0.1654957381714162
$ws.Range("M10").Value = 479.8918966666667
$ws.Range("N10").Value = 1439.67569
$ws.Range("O10").Value = 0.998322229344922
$ws.Range("P10").Value = 0.998322229344922
$ws.Range("Q10").Value = 6966.297449961136
$ws.Range("R10").Value = 62696.67704965023
$ws.Range("S10").Value = 0.1652180742783717
$ws.Range("T10").Value = 0.1652180742783717

# Row 11
$ws.Range("G11").Value = 1.46316
$ws.Range("H11").Value = 4.389480000000001
$ws.Range("I11").Value = 0.01668092142356404
$ws.Range("J11").Value = 0.01668092142356404
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.4526816666666666
$ws.Range("N11").Value = 1.358045
$ws.Range("O11").Value = 0.0009417166111561725
$ws.Range("P11").Value = 0.0009417166111561725
$ws.Range("Q11").Value = 0.6623457074000001
$ws.Range("R11").Value = 5.961111366600001
$ws.Range("S11").Value = 0.00001570870079396112
$ws.Range("T11").Value = 0.00001570870079396112

# Row 12
$ws.Range("G12").Value = 1.46316
$ws.Range("H12").Value = 4.389480000000001
$ws.Range("I12").Value = 0.01668092142356404
$ws.Range("J12").Value = 0.01668092142356404
$ws.Range("M12").Value = 0.35382
$ws.Range("N12").Value = 1.06146
$ws.Range("O12").Value = 0.0007360540439218367
$ws.Range("P12").Value = 0.0007360540439218367
$ws.Range("Q12").Value = 0.5176952712
$ws.Range("R12").Value = 4.6592574408
$ws.Range("S12").Value = 0.00001227805967015671
$ws.Range("T12").Value = 0.00001227805967015671

# Row 13
$ws.Range("G13").Value = 1.46316
$ws.Range("H13").Value = 4.389480000000001
$ws.Range("I13").Value = 0.01668092142356404
$ws.Range("J13").Value = 0.01668092142356404
$ws.Range("M13").Value = 479.8918966666667
$ws.Range("N13").Value = 1439.67569
$ws.Range("O13").Value = 0.998322229344922
$ws.Range("P13").Value = 0.998322229344922
$ws.Range("Q13").Value = 702.1586275268002
$ws.Range("R13").Value = 6319.427647741201
$ws.Range("S13").Value = 0.01665293466309992
$ws.Range("T13").Value = 0.01665293466309992
